$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2739
$ws.Range("F4").Value = 1069
$ws.Range("F5").Value = 19736
$ws.Range("F6").Value = 77
$ws.Range("F7").Value = 2268
$ws.Range("F8").Value = 750
$ws.Range("F9").Value = 613
$ws.Range("F10").Value = 438
$ws.Range("F11").Value = 693
$ws.Range("F12").Value = 238
$ws.Range("G13").Value = "已售罄"
$ws.Range("F15").Value = 369
$ws.Range("F16").Value = 73
$ws.Range("F17").Value = 263
$ws.Range("F19").Value = 194
$ws.Range("F22").Value = 101
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 197
$ws.Range("F4").Value = 12
$ws.Range("F7").Value = 288
$ws.Range("F12").Value = 11
$ws.Range("F13").Value = 87
$ws.Range("F15").Value = 75
$ws.Range("F20").Value = 18
$ws.Range("F22").Value = 22
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6004
$ws.Range("F3").Value = 646
$ws.Range("F4").Value = 592
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6004
$ws.Range("F3").Value = 646
$ws.Range("F4").Value = 592
$ws.Range("F5").Value = 197
$ws.Range("F8").Value = 2739
$ws.Range("F9").Value = 1069
$ws.Range("F10").Value = 19736
$ws.Range("F11").Value = 12
$ws.Range("F13").Value = 77
$ws.Range("F15").Value = 288
$ws.Range("F16").Value = 2268
$ws.Range("F17").Value = 750
$ws.Range("F19").Value = 613
$ws.Range("F20").Value = 438
$ws.Range("F21").Value = 693
$ws.Range("F22").Value = 238
$ws.Range("G23").Value = "已售罄"
$ws.Range("F28").Value = 369
$ws.Range("F29").Value = 73
$ws.Range("F30").Value = 11
$ws.Range("F31").Value = 263
$ws.Range("F32").Value = 87
$ws.Range("F35").Value = 194
$ws.Range("F36").Value = 75
$ws.Range("F43").Value = 18
$ws.Range("F45").Value = 22
$ws.Range("F47").Value = 101
